$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-17 Monday" "2025-11-18 Tuesday"

Replace-Text "45÷9=5, 0" "85÷4=21, 1"
Replace-Text "60÷4=15, 0" "33÷7=4, 5"
Replace-Text "45÷5=9, 0" "75÷6=12, 3"
Replace-Text "58÷4=14, 2" "90÷9=10, 0"
Replace-Text "48÷2=24, 0" "87÷8=10, 7"

Replace-Text "26÷9=2, 8" "47÷8=5, 7"
Replace-Text "73÷5=14, 3" "43÷3=14, 1"
Replace-Text "84÷3=28, 0" "97÷6=16, 1"
Replace-Text "30÷9=3, 3" "77÷4=19, 1"
Replace-Text "43÷4=10, 3" "77÷5=15, 2"

Replace-Text "67÷4=16, 3" "50÷2=25, 0"
Replace-Text "33÷2=16, 1" "97÷9=10, 7"
Replace-Text "78÷6=13, 0" "89÷8=11, 1"
Replace-Text "64÷6=10, 4" "91÷7=13, 0"
Replace-Text "93÷6=15, 3" "48÷4=12, 0"

Replace-Text "88÷9=9, 7" "17÷6=2, 5"
Replace-Text "68÷6=11, 2" "12÷6=2, 0"
Replace-Text "58÷3=19, 1" "59÷4=14, 3"
Replace-Text "76÷8=9, 4" "36÷4=9, 0"
Replace-Text "35÷9=3, 8" "43÷3=14, 1"

Replace-Text "10÷9=1, 1" "26÷7=3, 5"
Replace-Text "86÷5=17, 1" "97÷7=13, 6"
Replace-Text "83÷8=10, 3" "32÷8=4, 0"
Replace-Text "18÷2=9, 0" "77÷9=8, 5"
Replace-Text "90÷5=18, 0" "82÷8=10, 2"
